$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string values in columns A (emoji markers) and B (color labels)
# per commit: refactor the 4 marker emoji + rename "noir" -> "bleu"
$ws.Range("A2:A3").Value = "📕"
$ws.Range("A4:A11").Value = "📘"
$ws.Range("A12").Value = "📙"
$ws.Range("A13:A34").Value = "📘"
$ws.Range("A35").Value = "📗"
$ws.Range("A36:A49").Value = "📘"
$ws.Range("A50:A51").Value = "📗"
$ws.Range("A52:A58").Value = "📘"
$ws.Range("A59").Value = "📗"
$ws.Range("A60:A65").Value = "📘"
$ws.Range("B4:B11").Value = "bleu"
$ws.Range("B13:B34").Value = "bleu"
$ws.Range("B36:B49").Value = "bleu"
$ws.Range("B52:B58").Value = "bleu"
$ws.Range("B60:B65").Value = "bleu"
